# Add a new "2022-Q3" quarterly sheet to the workbook:
#  - insert it right after "总计", shifting "2021-Q4" / "2021-Q3" / "2020-Q4" back
#  - add a matching summary row on "总计"
#  - populate "2022-Q3" with its fund-holding detail rows

$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item(1)     # 总计 (totals sheet)
$q4_2021 = $wb.Worksheets.Item(2)  # 2021-Q4 (used as a formatting template)

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3
# ---------------------------------------------------------------------------
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()

# Give A2 the same (bold / bordered / centered) look as the other A-column cells
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 10
$zj.Range("D2").Value = 4.96

# The "A" column is a 0-based row index; renumber the rows that shifted down
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet by duplicating "2021-Q4" (so that it
#    automatically inherits the correct layout / styles), then overwrite
#    its contents with the 2022-Q3 fund-holding data.
# ---------------------------------------------------------------------------
$q4_2021.Copy($null, $zj)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template only has 3 rows (1 header + 2 data); we need 11 (1 header + 10 data)
$newSheet.Range("A4:A11").EntireRow.Insert()
$newSheet.Range("A3").Copy()
$newSheet.Range("A4:A11").PasteSpecial(-4122)

# Force columns B:G to be stored as text (several values look numeric, e.g.
# "83.65", but must stay text) without leaving a lingering number format.
$newSheet.Range("B2:G11").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002340"
$newSheet.Range("C2").Value = "富国价值优势混合"
$newSheet.Range("D2").Value = "83.65"
$newSheet.Range("E2").Value = "91.92"
$newSheet.Range("F2").Value = "2.63"
$newSheet.Range("G2").Value = "2.2000"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "004674"
$newSheet.Range("C3").Value = "富国新机遇灵活配置混合A"
$newSheet.Range("D3").Value = "36.29"
$newSheet.Range("E3").Value = "93.03"
$newSheet.Range("F3").Value = "2.42"
$newSheet.Range("G3").Value = "0.8782"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "519732"
$newSheet.Range("C4").Value = "交银定期支付双息平衡混合"
$newSheet.Range("D4").Value = "40.09"
$newSheet.Range("E4").Value = "66.00"
$newSheet.Range("F4").Value = "2.14"
$newSheet.Range("G4").Value = "0.8579"
$newSheet.Range("H4").Value = 5

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "012578"
$newSheet.Range("C5").Value = "富国红利混合A"
$newSheet.Range("D5").Value = "13.58"
$newSheet.Range("E5").Value = "93.23"
$newSheet.Range("F5").Value = "3.30"
$newSheet.Range("G5").Value = "0.4481"
$newSheet.Range("H5").Value = 3

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "011481"
$newSheet.Range("C6").Value = "广发瑞锦一年定开混合"
$newSheet.Range("D6").Value = "2.62"
$newSheet.Range("E6").Value = "89.31"
$newSheet.Range("F6").Value = "5.37"
$newSheet.Range("G6").Value = "0.1407"
$newSheet.Range("H6").Value = 6

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "012579"
$newSheet.Range("C7").Value = "富国红利混合C"
$newSheet.Range("D7").Value = "3.57"
$newSheet.Range("E7").Value = "93.23"
$newSheet.Range("F7").Value = "3.30"
$newSheet.Range("G7").Value = "0.1178"
$newSheet.Range("H7").Value = 3

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "004604"
$newSheet.Range("C8").Value = "富国新活力灵活配置混合A"
$newSheet.Range("D8").Value = "2.92"
$newSheet.Range("E8").Value = "93.16"
$newSheet.Range("F8").Value = "3.52"
$newSheet.Range("G8").Value = "0.1028"
$newSheet.Range("H8").Value = 7

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "004675"
$newSheet.Range("C9").Value = "富国新机遇灵活配置混合C"
$newSheet.Range("D9").Value = "3.93"
$newSheet.Range("E9").Value = "93.03"
$newSheet.Range("F9").Value = "2.42"
$newSheet.Range("G9").Value = "0.0951"
$newSheet.Range("H9").Value = 6

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "004605"
$newSheet.Range("C10").Value = "富国新活力灵活配置混合C"
$newSheet.Range("D10").Value = "2.17"
$newSheet.Range("E10").Value = "93.16"
$newSheet.Range("F10").Value = "3.52"
$newSheet.Range("G10").Value = "0.0764"
$newSheet.Range("H10").Value = 7

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "620001"
$newSheet.Range("C11").Value = "金元顺安宝石动力混合"
$newSheet.Range("D11").Value = "1.01"
$newSheet.Range("E11").Value = "40.12"
$newSheet.Range("F11").Value = "4.41"
$newSheet.Range("G11").Value = "0.0445"
$newSheet.Range("H11").Value = 3

# Drop the temporary text-number-format so these cells end up with no
# explicit style (matching the rest of the data rows), while keeping the
# values stored as text.
$newSheet.Range("B2:G11").ClearFormats()

# Header row labels (already correct from the template copy, re-assert to be safe)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

Write-Host "2022-Q3 sheet added."
